# Tutorial 6 solution update: switch the date separator in column A from
# "/" to "-" and refresh the attendance tally columns (D..H) for the rows
# whose counts changed.
#
# The date strings are forced to stay as TEXT (NumberFormat "@") before the
# write, then ClearFormats() removes the temporary formatting so the cell
# ends up with no explicit style - matching the original inlineStr cells
# which never carried a style index. Without this, Excel's smart-parsing
# would silently turn unambiguous values like "01-08-2022" into a real
# date serial (since day <= 12 could be read as a month).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DateText([string]$cellRef, [string]$text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

# Row 3 : 28/07/2022 -> 28-07-2022
Set-DateText "A3" "28-07-2022"
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

# Row 4 : 01/08/2022 -> 01-08-2022
Set-DateText "A4" "01-08-2022"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0

# Row 5 : 04/08/2022 -> 04-08-2022
Set-DateText "A5" "04-08-2022"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0

# Row 6 : 08/08/2022 -> 08-08-2022
Set-DateText "A6" "08-08-2022"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("H6").Value = 0

# Row 7 : 11/08/2022 -> 11-08-2022 (date only)
Set-DateText "A7" "11-08-2022"

# Row 8 : 15/08/2022 -> 15-08-2022 (date only)
Set-DateText "A8" "15-08-2022"

# Row 9 : 18/08/2022 -> 18-08-2022 (date only)
Set-DateText "A9" "18-08-2022"

# Row 10 : 22/08/2022 -> 22-08-2022
Set-DateText "A10" "22-08-2022"
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("H10").Value = 0

# Row 11 : 25/08/2022 -> 25-08-2022
Set-DateText "A11" "25-08-2022"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("H11").Value = 0

# Row 12 : 29/08/2022 -> 29-08-2022
Set-DateText "A12" "29-08-2022"
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("H12").Value = 0

# Row 13 : 01/09/2022 -> 01-09-2022
Set-DateText "A13" "01-09-2022"
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("H13").Value = 0

# Row 14 : 05/09/2022 -> 05-09-2022 (date only)
Set-DateText "A14" "05-09-2022"

# Row 15 : 08/09/2022 -> 08-09-2022 (date only)
Set-DateText "A15" "08-09-2022"

# Row 16 : 12/09/2022 -> 12-09-2022
Set-DateText "A16" "12-09-2022"
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 1
$ws.Range("H16").Value = 0

# Row 17 : 15/09/2022 -> 15-09-2022 (date only)
Set-DateText "A17" "15-09-2022"

# Row 18 : 19/09/2022 -> 19-09-2022 (date only)
Set-DateText "A18" "19-09-2022"

# Row 19 : 22/09/2022 -> 22-09-2022 (date only)
Set-DateText "A19" "22-09-2022"

# Row 20 : 26/09/2022 -> 26-09-2022
Set-DateText "A20" "26-09-2022"
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 1
$ws.Range("H20").Value = 0

# Row 21 : 29/09/2022 -> 29-09-2022 (date only)
Set-DateText "A21" "29-09-2022"
